$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "2025-11-26 00:44:49"
$ws.Range("C25").Value = "A873-150925-CHK-Y06"
$ws.Range("D25").Value = 993
$ws.Range("E25").Value = 993
$ws.Range("F25").Value = 951
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 25
$ws.Range("I25").Value = 17
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 813
$ws.Range("L25").Value = $false
$ws.Range("M25").Value = 42
$ws.Range("N25").Value = 4.23
$ws.Range("O25").Value = 1.96
